$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H111").Value = 1614.8334
$ws.Range("J111").Value = 2047.25
$ws.Range("L111").Value = 6141.75
$ws.Range("N111").Value = -12275.75
$ws.Range("H112").Value = 1302.6461
$ws.Range("J112").Value = 1302.6461
$ws.Range("L112").Value = 3907.9383
$ws.Range("N112").Value = -6123.9383
$ws.Range("H129").Value = 1853.4147
$ws.Range("J129").Value = 2167.4412
$ws.Range("L129").Value = 6502.323600000001
$ws.Range("N129").Value = -16502.3236
$ws.Range("H135").Value = 1096
$ws.Range("I135").Value = 733.4286
$ws.Range("J135").Value = 2365
$ws.Range("K135").Value = 6600.8574
$ws.Range("L135").Value = 21285
$ws.Range("M135").Value = -4065.8574
$ws.Range("N135").Value = -26355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2437.2407
$ws.Range("I132").Value = 1554.7241
$ws.Range("J132").Value = 3460.96
$ws.Range("K132").Value = 4664.1723
$ws.Range("L132").Value = 10382.88
$ws.Range("M132").Value = -2134.1723
$ws.Range("N132").Value = -15442.88
$ws.Range("H137").Value = 47580
$ws.Range("J137").Value = 47580
$ws.Range("L137").Value = 47580
$ws.Range("N137").Value = -57780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 50282
$ws.Range("J132").Value = 50282
$ws.Range("L132").Value = 50282
$ws.Range("N132").Value = -60402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3192.7144
$ws.Range("I134").Value = 1069.8
$ws.Range("J134").Value = 8500
$ws.Range("K134").Value = 3209.4
$ws.Range("L134").Value = 25500
$ws.Range("M134").Value = -674.3999999999996
$ws.Range("N134").Value = -30570
$ws.Range("H137").Value = 44540
$ws.Range("J137").Value = 44540
$ws.Range("L137").Value = 44540
$ws.Range("N137").Value = -54740
$ws.Range("H138").Value = 44450
$ws.Range("J138").Value = 44450
$ws.Range("L138").Value = 44450
$ws.Range("N138").Value = -54730
$ws.Range("H140").Value = 120851.25
$ws.Range("J140").Value = 120851.25
$ws.Range("L140").Value = 120851.25
$ws.Range("N140").Value = -131211.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 120
$ws.Range("I7").Value = 120
$ws.Range("K7").Value = 360
$ws.Range("M7").Value = -248
$ws.Range("H9").Value = 1145140
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1145140
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 3435420
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -3435868
$ws.Range("H63").Value = 3391.8096
$ws.Range("I63").Value = 2801
$ws.Range("J63").Value = 4573.4287
$ws.Range("K63").Value = 8403
$ws.Range("L63").Value = 13720.2861
$ws.Range("M63").Value = -7654
$ws.Range("N63").Value = -15218.2861
$ws.Range("H64").Value = 6984.857
$ws.Range("J64").Value = 8738.799999999999
$ws.Range("L64").Value = 26216.4
$ws.Range("N64").Value = -26756.4
$ws.Range("H66").Value = 3391.8096
$ws.Range("I66").Value = 2801
$ws.Range("J66").Value = 4573.4287
$ws.Range("K66").Value = 25209
$ws.Range("L66").Value = 41160.85830000001
$ws.Range("M66").Value = -21465
$ws.Range("N66").Value = -48648.85830000001
$ws.Range("H67").Value = 6984.857
$ws.Range("J67").Value = 8738.799999999999
$ws.Range("L67").Value = 26216.4
$ws.Range("N67").Value = -28088.4
$ws.Range("H68").Value = 1127.0377
$ws.Range("J68").Value = 1247.75
$ws.Range("L68").Value = 3743.25
$ws.Range("N68").Value = -5365.25
$ws.Range("H69").Value = 154961.53
$ws.Range("I69").Value = 737.5
$ws.Range("J69").Value = 401720
$ws.Range("K69").Value = 2212.5
$ws.Range("L69").Value = 1205160
$ws.Range("M69").Value = -1401.5
$ws.Range("N69").Value = -1206782
$ws.Range("H71").Value = 1127.0377
$ws.Range("J71").Value = 1247.75
$ws.Range("L71").Value = 11229.75
$ws.Range("N71").Value = -19341.75
$ws.Range("H72").Value = 154961.53
$ws.Range("I72").Value = 737.5
$ws.Range("J72").Value = 401720
$ws.Range("K72").Value = 6637.5
$ws.Range("L72").Value = 3615480
$ws.Range("M72").Value = -2581.5
$ws.Range("N72").Value = -3623592
$ws.Range("H101").Value = 5000
$ws.Range("J101").Value = 5000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -19868
$ws.Range("H106").Value = 3414.8462
$ws.Range("J106").Value = 3414.8462
$ws.Range("L106").Value = 10244.5386
$ws.Range("N106").Value = -12136.5386
$ws.Range("H113").Value = 571.88635
$ws.Range("I113").Value = 481.10526
$ws.Range("J113").Value = 640.88
$ws.Range("K113").Value = 1443.31578
$ws.Range("L113").Value = 1922.64
$ws.Range("M113").Value = 726.6842200000001
$ws.Range("N113").Value = -6262.639999999999
$ws.Range("H131").Value = 765.63635
$ws.Range("I131").Value = 433.33334
$ws.Range("J131").Value = 798.86664
$ws.Range("K131").Value = 1300.00002
$ws.Range("L131").Value = 2396.59992
$ws.Range("M131").Value = 3739.99998
$ws.Range("N131").Value = -12476.59992
$ws.Range("H132").Value = 1917.2941
$ws.Range("I132").Value = 670.2
$ws.Range("J132").Value = 2436.9167
$ws.Range("K132").Value = 6031.8
$ws.Range("L132").Value = 21932.2503
$ws.Range("M132").Value = -3501.8
$ws.Range("N132").Value = -26992.2503
$ws.Range("H136").Value = 3611.3333
$ws.Range("I136").Value = 2446.25
$ws.Range("J136").Value = 4942.857
$ws.Range("K136").Value = 7338.75
$ws.Range("L136").Value = 14828.571
$ws.Range("M136").Value = -2238.75
$ws.Range("N136").Value = -25028.571
$ws.Range("H137").Value = 3608.7742
$ws.Range("I137").Value = 2646.6667
$ws.Range("J137").Value = 4940.923
$ws.Range("K137").Value = 7940.000100000001
$ws.Range("L137").Value = 14822.769
$ws.Range("M137").Value = -2840.000100000001
$ws.Range("N137").Value = -25022.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5701.9824
$ws.Range("I70").Value = 5480.674
$ws.Range("J70").Value = 6627.4546
$ws.Range("K70").Value = 5480.674
$ws.Range("L70").Value = 6627.4546
$ws.Range("M70").Value = -5210.674
$ws.Range("N70").Value = -7167.4546
$ws.Range("H73").Value = 5701.9824
$ws.Range("I73").Value = 5480.674
$ws.Range("J73").Value = 6627.4546
$ws.Range("K73").Value = 5480.674
$ws.Range("L73").Value = 6627.4546
$ws.Range("M73").Value = -4544.674
$ws.Range("N73").Value = -8499.454600000001
$ws.Range("H122").Value = 6236.0454
$ws.Range("I122").Value = 3126.6365
$ws.Range("J122").Value = 9345.454
$ws.Range("K122").Value = 9379.9095
$ws.Range("L122").Value = 28036.362
$ws.Range("M122").Value = -6929.9095
$ws.Range("N122").Value = -32936.362
$ws.Range("H132").Value = 4863.3184
$ws.Range("I132").Value = 4050.111
$ws.Range("J132").Value = 5426.3076
$ws.Range("K132").Value = 12150.333
$ws.Range("L132").Value = 16278.9228
$ws.Range("M132").Value = -9620.332999999999
$ws.Range("N132").Value = -21338.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 364.4
$ws.Range("I55").Value = 392.33334
$ws.Range("J55").Value = 322.5
$ws.Range("K55").Value = 392.33334
$ws.Range("L55").Value = 322.5
$ws.Range("M55").Value = -219.33334
$ws.Range("N55").Value = -668.5
$ws.Range("H69").Value = 39550
$ws.Range("J69").Value = 39550
$ws.Range("L69").Value = 39550
$ws.Range("N69").Value = -41172
$ws.Range("H72").Value = 39550
$ws.Range("J72").Value = 39550
$ws.Range("L72").Value = 118650
$ws.Range("N72").Value = -126762
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H122").Value = 4283.3335
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 4358.8237
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 13076.4711
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -17976.4711
$ws.Range("H132").Value = 3947.6191
$ws.Range("I132").Value = 2394.4443
$ws.Range("J132").Value = 13266.667
$ws.Range("K132").Value = 7183.3329
$ws.Range("L132").Value = 39800.001
$ws.Range("M132").Value = -4653.3329
$ws.Range("N132").Value = -44860.001
